# ============================================================
# edit.ps1 -- apply tense/wording revisions to IMPACT_OF_AUTO
# ============================================================
$d = $word.ActiveDocument

function Replace-One([string]$find, [string]$replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 1)
    if (-not $ok) { Write-Output "MISSING: $find" }
}

Replace-One "ol that is being deployed in an" "ol that was deployed in an"
Replace-One " introductory Computer Science course. The tool is a web site that accepts student submissions for assignments, and automatically tests them for correctness. The students receive feedback moments after submission. They may use that feedback to improve their work, and submit again. We compar" " introductory Computer Science course. The tool was a web site that accepted student submissions for assignments, and automatically tested them for correctness. The students received feedback moments after submission. They could use that feedback to improve their work, and submitted again. We compar"
Replace-One "We introduce a new" "We introduced a new"
Replace-One "n introductory computing course, as well as an assessment of the impact of that system on teaching. This system simplifies the logistics of the grading process and can help students by providing near-immediate" "n introductory computing course, and assessed the impact of that system on teaching. This system simplified the logistics of the grading process and helped students by providing near-immediate"
Replace-One "The assignments in this course are programming exercises written in C. The course typically has 40 to 60 i" "The assignments in this course were programming exercises written in C. The course typically has 40 to 60 i"
Replace-One " It is important to get feedback to the students as soon as possible, and this pace (coupled with the large enrollments of 50 to 60 students per section) creates a large workload for the instructor and the graders." " It is important to get feedback to the students as soon as possible [0], and this pace (coupled with the large enrollments of 50 to 60 students per section) creates a large workload for the instructor and the graders."
Replace-One " a pedagogical problem. Presently, a student only sees feedback of " " a pedagogical problem. Before this system, a student only saw feedback of "
Replace-One " is detrimental to students, inhibiting their ability for informed iteration." " was detrimental to students, inhibiting their ability for informed iteration."
Replace-One " If we can provide students with immediate feedback on their work, we can solve the pedagogical problem and use staff time more efficiently." " "
Replace-One " was initially developed to support the teaching of a “flipped” course, where students watch video lectures online before class to prepare for classroom questions and discussion. Bottlenose also includes online submission and grading of programming assignments, which turns out to be a useful piece of functionality even for traditional courses. " " was initially developed to support the teaching of a “flipped” course, where students watch video lectures online before class to prepare for classroom questions and discussion. It also included online submission and grading of programming assignments, which were useful functions for traditional courses, as is examined in this paper."
Replace-One "Bottlenose is built" "The system was built"
Replace-One "work. This framework has allowed the application to be built r" "work. This framework allowed the application to be built r"
Replace-One "in automated testing infrastructure which will help the application stay high quality and maintainable as it grows. The application follows standard Rails conventions. A PostgreSQL2 database is used to store most application state, although student submissions are stored on the file system. " "in automated testing infrastructure which has helped the application stay high quality and maintainable as it grew. The application followed standard Rails conventions. A PostgreSQL2 database was used to store most application state, although student submissions were stored on the file system. "
Replace-One "A simple process for online submission of assignments is provided. Students are emailed an authentication link that brings them to their list of assignments and identifies them to the application. Assignments are " "A simple process for online submission of assignments was provided. Students were emailed authentication links that brought them to their list of assignments, and identified the students to the application. Assignments were "
Replace-One ") are supported. The automated grading process begins immediately when an assignment is submitted, giving students feedback within a few seconds. Students may attempt submissions multiple times. " ") were supported. The automated grading process began immediately when an assignment was submitted, giving students feedback within a few seconds. Students could attempt submissions multiple times. "
Replace-One "In order to automatically grade student programs, submissions are compiled and run on the server. Allowing students to run arbitrary code on the server is clearly a potential security issue. Bottlenose uses a sandbox mechanism to prevent student programs from causing trouble. Five major techniques are used to isolate student programs from the rest of the system: " "In order to automatically grade student programs, submissions were compiled and run on the server. Allowing students to run arbitrary code on the server is clearly a potential security issue [4]. The system used a sandbox mechanism to prevent student programs from causing trouble. Five major techniques are used to isolate student programs from the rest of the system: "
Replace-One " - Each student program is run under a separate system user with minimal Unix permissions. " " - Each student program was run under a separate system user with minimal Unix permissions. "
Replace-One " - Student programs can only access specific, white" " - Student programs could only access specific, white"
Replace-One " - The “setrlimit” system call is used to set limits on the use of a variety of resources, including RAM, child processes, and created file size. " " - The “setrlimit” system call was used to set limits on the use of a variety of resources, including RAM, child processes, and created file size. "
Replace-One " - Each program is executed in a separate “tmpfs” file" " - Each program was executed in a separate “tmpfs” file"
Replace-One " which ceases to exist when the grading process finishes. " " which ceased to exist when the grading process finished. "
Replace-One " - A grading process is terminated if it lasts more than five minutes. " " - A grading process was terminated if it lasted more than five minutes. "
Replace-One "This sandbox mechanism does not " "This sandbox mechanism did not "
Replace-One "vulnerability. It does perform adequately " "vulnerability. It did perform adequately "
Replace-One "platform allows for any" "platform allowed for any"

# ------------------------------------------------------------
# Paragraph-level formatting: four paragraphs gain a first-line
# (firstLine=720 twips = 36 pt) indent now that they start a
# new paragraph of body text.
# ------------------------------------------------------------
function Set-FirstLineIndent-ByPrefix([string]$prefix) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text
        if ($t.Length -ge $prefix.Length -and $t.Substring(0, $prefix.Length) -eq $prefix) {
            $p.Range.ParagraphFormat.FirstLineIndent = 36
            return $true
        }
    }
    Write-Output "MISSING PARA: $prefix"
    return $false
}

[void](Set-FirstLineIndent-ByPrefix "The assignments in this course were")
[void](Set-FirstLineIndent-ByPrefix "The system was built")
[void](Set-FirstLineIndent-ByPrefix "A simple process for online submission")
[void](Set-FirstLineIndent-ByPrefix "In order to automatically grade student")
